# "Generate Report for Handback"
#
# The handback files for both localized-file rows have now been produced,
# so the status report is refreshed:
#   - Overview sheet: the summary status for each language goes from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - Per-language sheets (zh-cn, de-de): the "Latest Target File" column
#     gets a hyperlink to the source .md (same link as column A), the
#     "Latest Handback File" column is filled in with the generated
#     handback xliff file name, and "Latest Handback DateTime" records
#     when each handback was produced.

$wb = $excel.ActiveWorkbook

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f93060abcf824b8ce251879fcc79beecf6babd8d/e2e/7aaef43c-e049-47f7-89ff-982a0b7beab2.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f93060abcf824b8ce251879fcc79beecf6babd8d/e2e/efdf98a9-f36d-40f7-9be4-e2885403b883.md"
$mdName1 = "7aaef43c-e049-47f7-89ff-982a0b7beab2.md"
$mdName2 = "efdf98a9-f36d-40f7-9be4-e2885403b883.md"

# ---------------------------------------------------------------------------
# Overview sheet - update the per-language status summary now that both
# languages have been handed back in sync with the English source.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------------------
# zh-cn sheet - fill in the target file link, the handback xliff name, and
# the handback timestamp for both rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdName1)
$zhcn.Range("J2").Value = "7aaef43c-e049-47f7-89ff-982a0b7beab2.17d2ec0a5d6b6244811f8dbb1b324158e60ea921.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 05:09:53"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdName2)
$zhcn.Range("J3").Value = "efdf98a9-f36d-40f7-9be4-e2885403b883.210b98d1b630675c9c851c7ce938e40864dd7e64.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 05:09:53"

$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(9).AutoFit()
$zhcn.Columns.Item(10).AutoFit()

# ---------------------------------------------------------------------------
# de-de sheet - same refresh, with its own (later) handback timestamp.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdName1)
$dede.Range("J2").Value = "7aaef43c-e049-47f7-89ff-982a0b7beab2.17d2ec0a5d6b6244811f8dbb1b324158e60ea921.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 05:10:20"

$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdName2)
$dede.Range("J3").Value = "efdf98a9-f36d-40f7-9be4-e2885403b883.210b98d1b630675c9c851c7ce938e40864dd7e64.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 05:10:20"

$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(9).AutoFit()
$dede.Columns.Item(10).AutoFit()

Write-Output "Handback report generated."
